$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("invalid")

# Shift the existing action_taken / app_recipient columns (H:I) two columns
# to the right (J:K) to make room for the new credit_purpose /
# credit_purpose_ff columns.
$ws.Range("H1:I11").Copy()
$ws.Range("J1").PasteSpecial()
$ws.Range("H1:I11").Clear()

# New column headers.
$ws.Range("H1").Value = "credit_purpose"
$ws.Range("I1").Value = "credit_purpose_ff"

# New column data.
$ws.Range("H2").Value = "1;2;3"
$ws.Range("I2").Value = "abc;def"

$ws.Range("H3").Value = 988
$ws.Range("I3").Value = "abc;def"

$ws.Range("H4").Value = "1;2;3;4;5;6;7;8"
$ws.Range("I4").Value = "abc;def"

$ws.Range("H5").Value = "10000;1200;1;2"
$ws.Range("I5").Value = "abc;def"

$ws.Range("H6").Value = "1;999"
$ws.Range("I6").Value = "abc;def"

$ws.Range("H7").Value = "1;2;3;4;5;6"

$ws.Range("H8").Value = 977

$ws.Range("H9").Value = 977
$ws.Range("I9").Value = "abc;def"

$ws.Range("H10").Value = "977;1"
$ws.Range("I10").Value = "abc;def"

$ws.Range("H11").Value = 999
$ws.Range("I11").Value = "abc;def"

# The new credit_purpose column is a bit narrower than the rest.
$ws.Columns.Item(8).ColumnWidth = 13.67

# Rows 4 and 5 now wrap onto a second line.
$ws.Rows.Item(4).RowHeight = 34
$ws.Rows.Item(5).RowHeight = 34

# Restore the cursor to where the author left it.
$ws.Range("H7").Select()
